$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1826.9333
$ws.Range("I2").Value = 211.75
$ws.Range("J2").Value = 3672.8572
$ws.Range("K2").Value = 211.75
$ws.Range("L2").Value = 3672.8572
$ws.Range("M2").Value = -98.75
$ws.Range("N2").Value = -3898.8572
$ws.Range("H74").Value = 102387.625
$ws.Range("I74").Value = 139291.19
$ws.Range("J74").Value = 21199.8
$ws.Range("K74").Value = 139291.19
$ws.Range("L74").Value = 21199.8
$ws.Range("M74").Value = -138355.19
$ws.Range("N74").Value = -23071.8
$ws.Range("H77").Value = 102387.625
$ws.Range("I77").Value = 139291.19
$ws.Range("J77").Value = 21199.8
$ws.Range("K77").Value = 696455.95
$ws.Range("L77").Value = 105999
$ws.Range("M77").Value = -691775.95
$ws.Range("N77").Value = -115359
$ws.Range("H112").Value = 4425.6
$ws.Range("J112").Value = 4517.3335
$ws.Range("L112").Value = 13552.0005
$ws.Range("N112").Value = -15768.0005
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H131").Value = 56052.5
$ws.Range("I131").Value = 100000
$ws.Range("K131").Value = 300000
$ws.Range("M131").Value = -294960
$ws.Range("H137").Value = 1593.9
$ws.Range("I137").Value = 1563.5714
$ws.Range("K137").Value = 4690.7142
$ws.Range("M137").Value = -2140.7142
$ws.Range("H138").Value = 2878.1667
$ws.Range("I138").Value = 1726.4546
$ws.Range("K138").Value = 5179.3638
$ws.Range("M138").Value = -39.36380000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 2965292.2
$ws.Range("I32").Value = 2803904
$ws.Range("K32").Value = 2803904
$ws.Range("M32").Value = -2803617
$ws.Range("H45").Value = 4997.3335
$ws.Range("I45").Value = 5796.8
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 5796.8
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -5419.8
$ws.Range("N45").Value = -1754
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 3202.125
$ws.Range("I102").Value = 2945.2856
$ws.Range("K102").Value = 2945.2856
$ws.Range("M102").Value = -1323.2856
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 1995.7142
$ws.Range("I122").Value = 1661.6666
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4984.9998
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2534.9998
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H20").Value = 1428.3334
$ws.Range("I20").Value = 1428.3334
$ws.Range("K20").Value = 1428.3334
$ws.Range("M20").Value = -1181.3334
$ws.Range("H25").Value = 207
$ws.Range("I25").Value = 207
$ws.Range("K25").Value = 207
$ws.Range("M25").Value = 28
$ws.Range("H39").Value = 14998
$ws.Range("I39").Value = 14998
$ws.Range("K39").Value = 14998
$ws.Range("M39").Value = -14609
$ws.Range("H54").Value = 3472.5
$ws.Range("J54").Value = 4500
$ws.Range("L54").Value = 4500
$ws.Range("N54").Value = -5468
$ws.Range("H105").Value = 2666.6667
$ws.Range("I105").Value = 2666.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2666.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -919.6667000000002
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1986.2222
$ws.Range("I107").Value = 1986.2222
$ws.Range("K107").Value = 1986.2222
$ws.Range("M107").Value = -66.22219999999993

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2260.5625
$ws.Range("I31").Value = 2115
$ws.Range("J31").Value = 4444
$ws.Range("K31").Value = 2115
$ws.Range("L31").Value = 4444
$ws.Range("M31").Value = -1820
$ws.Range("N31").Value = -5034
$ws.Range("H34").Value = 2260.5625
$ws.Range("I34").Value = 2115
$ws.Range("J34").Value = 4444
$ws.Range("K34").Value = 2115
$ws.Range("L34").Value = 4444
$ws.Range("M34").Value = -1913
$ws.Range("N34").Value = -4848
$ws.Range("H58").Value = 1839.125
$ws.Range("I58").Value = 1625.7142
$ws.Range("K58").Value = 1625.7142
$ws.Range("M58").Value = -1422.7142
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H105").Value = 3062.5
$ws.Range("I105").Value = 2200
$ws.Range("K105").Value = 2200
$ws.Range("M105").Value = -453
$ws.Range("H107").Value = 1150.1666
$ws.Range("I107").Value = 997
$ws.Range("K107").Value = 997
$ws.Range("M107").Value = 923
$ws.Range("H136").Value = 1839.125
$ws.Range("I136").Value = 1625.7142
$ws.Range("K136").Value = 4877.142599999999
$ws.Range("M136").Value = -2327.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 799
$ws.Range("J117").Value = 699.5
$ws.Range("L117").Value = 2098.5
$ws.Range("N117").Value = -8982.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 4353.7144
$ws.Range("I132").Value = 4197.2
$ws.Range("K132").Value = 37774.8
$ws.Range("M132").Value = -35244.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 25058
$ws.Range("J94").Value = 27053.727
$ws.Range("L94").Value = 27053.727
$ws.Range("N94").Value = -28405.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8795.143
$ws.Range("I7").Value = 10347
$ws.Range("K7").Value = 10347
$ws.Range("M7").Value = -10235
$ws.Range("H46").Value = 1920.5
$ws.Range("I46").Value = 1099.3334
$ws.Range("J46").Value = 2194.2222
$ws.Range("K46").Value = 1099.3334
$ws.Range("L46").Value = 2194.2222
$ws.Range("M46").Value = -911.3334
$ws.Range("N46").Value = -2570.2222
$ws.Range("H55").Value = 1230
$ws.Range("I55").Value = 654.6
$ws.Range("J55").Value = 1949.25
$ws.Range("K55").Value = 654.6
$ws.Range("L55").Value = 1949.25
$ws.Range("M55").Value = -481.6
$ws.Range("N55").Value = -2295.25
$ws.Range("H68").Value = 1191.6923
$ws.Range("J68").Value = 2666
$ws.Range("L68").Value = 2666
$ws.Range("N68").Value = -4164
$ws.Range("H71").Value = 1191.6923
$ws.Range("J71").Value = 2666
$ws.Range("L71").Value = 13330
$ws.Range("N71").Value = -20818
$ws.Range("H93").Value = 1419.2
$ws.Range("I93").Value = 1260.6
$ws.Range("K93").Value = 1260.6
$ws.Range("M93").Value = -12.59999999999991
$ws.Range("H126").Value = 8795.143
$ws.Range("I126").Value = 10347
$ws.Range("K126").Value = 31041
$ws.Range("M126").Value = -28571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13332.429
$ws.Range("I41").Value = 10839.667
$ws.Range("J41").Value = 15202
$ws.Range("K41").Value = 10839.667
$ws.Range("L41").Value = 15202
$ws.Range("M41").Value = -10449.667
$ws.Range("N41").Value = -15982
$ws.Range("H126").Value = 3163.889
$ws.Range("I126").Value = 1855.6
$ws.Range("K126").Value = 5566.799999999999
$ws.Range("M126").Value = -3096.799999999999
$ws.Range("H132").Value = 1356.7142
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 4249.0002
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -1719.0002
$ws.Range("N132").Value = -8057
